$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append rows 637-657: continuation of the NAV time series.
# The new block repeats the same underlying market data as the existing
# 2024-08-28 .. 2024-09-25 block (rows 616-636), with column J (NAV) carried
# forward by compounding on top of the prior last NAV value (row 636).

$ws.Range("A637").Value = "'2024-08-28"
$ws.Range("C637").Value = 1713.5
$ws.Range("D637").Value = 611.2000122070312
$ws.Range("E637").Value = 1138.300048828125
$ws.Range("F637").Value = 180.9600067138672
$ws.Range("G637").Value = 641.5499877929688
$ws.Range("H637").Value = 22112.82034301758
$ws.Range("I637").Value = 0
$ws.Range("I637").NumberFormat = "General"
$ws.Range("J637").Value = 223.292769200981

$ws.Range("A638").Value = "'2024-08-29"
$ws.Range("C638").Value = 1755.650024414062
$ws.Range("D638").Value = 603.6199951171875
$ws.Range("E638").Value = 1132.050048828125
$ws.Range("F638").Value = 179.9400024414062
$ws.Range("G638").Value = 644.2999877929688
$ws.Range("H638").Value = 22256.81030273438
$ws.Range("I638").Value = 0.006511605371146772
$ws.Range("I638").NumberFormat = "General"
$ws.Range("J638").Value = 224.7467635962483

$ws.Range("A639").Value = "'2024-08-30"
$ws.Range("C639").Value = 1783.050048828125
$ws.Range("D639").Value = 600.3599853515625
$ws.Range("E639").Value = 1127.900024414062
$ws.Range("F639").Value = 178.6199951171875
$ws.Range("G639").Value = 632.0499877929688
$ws.Range("H639").Value = 22286.13012695312
$ws.Range("I639").Value = 0.001317341695415713
$ws.Range("I639").NumberFormat = "General"
$ws.Range("J639").Value = 225.0428318788434

$ws.Range("A640").Value = "'2024-09-02"
$ws.Range("C640").Value = 1840.550048828125
$ws.Range("D640").Value = 608.5800170898438
$ws.Range("E640").Value = 1111.550048828125
$ws.Range("F640").Value = 177.5399932861328
$ws.Range("G640").Value = 670.2000122070312
$ws.Range("H640").Value = 22669.95040893555
$ws.Range("I640").Value = 0.01722238359894636
$ws.Range("I640").NumberFormat = "General"
$ws.Range("J640").Value = 228.9186058556541

$ws.Range("A641").Value = "'2024-09-03"
$ws.Range("C641").Value = 1865.599975585938
$ws.Range("D641").Value = 599.9400024414062
$ws.Range("E641").Value = 1114
$ws.Range("F641").Value = 178.4600067138672
$ws.Range("G641").Value = 659.0999755859375
$ws.Range("H641").Value = 22746.27993774414
$ws.Range("I641").Value = 0.003366991432787071
$ws.Range("I641").NumberFormat = "General"
$ws.Range("J641").Value = 229.6893728403757

$ws.Range("A642").Value = "'2024-09-04"
$ws.Range("C642").Value = 1871.900024414062
$ws.Range("D642").Value = 609
$ws.Range("E642").Value = 1127.900024414062
$ws.Range("F642").Value = 176.0200042724609
$ws.Range("G642").Value = 650.8499755859375
$ws.Range("H642").Value = 22782.94021606445
$ws.Range("I642").Value = 0.001611704349926693
$ws.Range("I642").NumberFormat = "General"
$ws.Range("J642").Value = 230.0595642017144

$ws.Range("A643").Value = "'2024-09-05"
$ws.Range("C643").Value = 1864.949951171875
$ws.Range("D643").Value = 602.1799926757812
$ws.Range("E643").Value = 1115.150024414062
$ws.Range("F643").Value = 173.4799957275391
$ws.Range("G643").Value = 643.8499755859375
$ws.Range("H643").Value = 22586.02963256836
$ws.Range("I643").Value = -0.008642896027846764
$ws.Range("I643").NumberFormat = "General"
$ws.Range("J643").Value = 228.0711833081073

$ws.Range("A644").Value = "'2024-09-06"
$ws.Range("C644").Value = 1857.150024414062
$ws.Range("D644").Value = 597.2999877929688
$ws.Range("E644").Value = 1100
$ws.Range("F644").Value = 169.8500061035156
$ws.Range("G644").Value = 634.7000122070312
$ws.Range("H644").Value = 22350.45025634766
$ws.Range("I644").Value = -0.01043031378480992
$ws.Range("I644").NumberFormat = "General"
$ws.Range("J644").Value = 225.6923293009308

$ws.Range("A645").Value = "'2024-09-09"
$ws.Range("C645").Value = 1860.449951171875
$ws.Range("D645").Value = 610.3400268554688
$ws.Range("E645").Value = 1104.150024414062
$ws.Range("F645").Value = 168.3300018310547
$ws.Range("G645").Value = 635.2000122070312
$ws.Range("H645").Value = 22404.27005004883
$ws.Range("I645").Value = 0.002407995950143633
$ws.Range("I645").NumberFormat = "General"
$ws.Range("J645").Value = 226.2357955158659

$ws.Range("A646").Value = "'2024-09-10"
$ws.Range("C646").Value = 1824.5
$ws.Range("D646").Value = 608
$ws.Range("E646").Value = 1113.199951171875
$ws.Range("F646").Value = 169.75
$ws.Range("G646").Value = 637.0499877929688
$ws.Range("H646").Value = 22289.99975585938
$ws.Range("I646").Value = -0.005100380147810443
$ws.Range("I646").NumberFormat = "General"
$ws.Range("J646").Value = 225.0819069556927

$ws.Range("A647").Value = "'2024-09-11"
$ws.Range("C647").Value = 1833.150024414062
$ws.Range("D647").Value = 627.6599731445312
$ws.Range("E647").Value = 1112.599975585938
$ws.Range("F647").Value = 165.8800048828125
$ws.Range("G647").Value = 627.2000122070312
$ws.Range("H647").Value = 22284.95007324219
$ws.Range("I647").Value = -0.0002265447587481507
$ws.Range("I647").NumberFormat = "General"
$ws.Range("J647").Value = 225.0309158293828

$ws.Range("A648").Value = "'2024-09-12"
$ws.Range("C648").Value = 1854.849975585938
$ws.Range("D648").Value = 645.5999755859375
$ws.Range("E648").Value = 1120.099975585938
$ws.Range("F648").Value = 167.0200042724609
$ws.Range("G648").Value = 651.0999755859375
$ws.Range("H648").Value = 22615.88967895508
$ws.Range("I648").Value = 0.01485036334500268
$ws.Range("I648").NumberFormat = "General"
$ws.Range("J648").Value = 228.3727066933079

$ws.Range("A649").Value = "'2024-09-13"
$ws.Range("C649").Value = 1894.449951171875
$ws.Range("D649").Value = 646.6500244140625
$ws.Range("E649").Value = 1118.550048828125
$ws.Range("F649").Value = 167.25
$ws.Range("G649").Value = 633.4500122070312
$ws.Range("H649").Value = 22746.35009765625
$ws.Range("I649").Value = 0.005768529142701387
$ws.Range("I649").NumberFormat = "General"
$ws.Range("J649").Value = 229.6900813072659

$ws.Range("A650").Value = "'2024-09-16"
$ws.Range("C650").Value = 1857.599975585938
$ws.Range("D650").Value = 621.0499877929688
$ws.Range("E650").Value = 1115.849975585938
$ws.Range("F650").Value = 163.9600067138672
$ws.Range("G650").Value = 665.9500122070312
$ws.Range("H650").Value = 22506.51992797852
$ws.Range("I650").Value = -0.01054367705799297
$ws.Range("I650").NumberFormat = "General"
$ws.Range("J650").Value = 227.2683032665379

$ws.Range("A651").Value = "'2024-09-17"
$ws.Range("C651").Value = 1848.699951171875
$ws.Range("D651").Value = 649.6500244140625
$ws.Range("E651").Value = 1110.949951171875
$ws.Range("F651").Value = 160.6000061035156
$ws.Range("G651").Value = 666.3499755859375
$ws.Range("H651").Value = 22484.49969482422
$ws.Range("I651").Value = -0.0009783935155129372
$ws.Range("I651").NumberFormat = "General"
$ws.Range("J651").Value = 227.0459454323403

$ws.Range("A652").Value = "'2024-09-18"
$ws.Range("C652").Value = 1888.199951171875
$ws.Range("D652").Value = 646.7000122070312
$ws.Range("E652").Value = 1079.949951171875
$ws.Range("F652").Value = 158.5599975585938
$ws.Range("G652").Value = 651.7000122070312
$ws.Range("H652").Value = 22442.71960449219
$ws.Range("I652").Value = -0.00185817300358472
$ws.Range("I652").NumberFormat = "General"
$ws.Range("J652").Value = 226.6240547859645

$ws.Range("A653").Value = "'2024-09-19"
$ws.Range("C653").Value = 1890.400024414062
$ws.Range("D653").Value = 652.1500244140625
$ws.Range("E653").Value = 1054.449951171875
$ws.Range("F653").Value = 155.25
$ws.Range("G653").Value = 649.5999755859375
$ws.Range("H653").Value = 22292.29992675781
$ws.Range("I653").Value = -0.006702381903139165
$ws.Range("I653").NumberFormat = "General"
$ws.Range("J653").Value = 225.1051338223511

$ws.Range("A654").Value = "'2024-09-20"
$ws.Range("C654").Value = 1916.800048828125
$ws.Range("D654").Value = 654.4500122070312
$ws.Range("E654").Value = 1054.599975585938
$ws.Range("F654").Value = 161.4299926757812
$ws.Range("G654").Value = 665.1500244140625
$ws.Range("H654").Value = 22632.26013183594
$ws.Range("I654").Value = 0.01525011803156592
$ws.Range("I654").NumberFormat = "General"
$ws.Range("J654").Value = 228.5380136826533

$ws.Range("A655").Value = "'2024-09-23"
$ws.Range("C655").Value = 1919.949951171875
$ws.Range("D655").Value = 654.0999755859375
$ws.Range("E655").Value = 1055.25
$ws.Range("F655").Value = 159.5599975585938
$ws.Range("G655").Value = 672
$ws.Range("H655").Value = 22635.46960449219
$ws.Range("I655").Value = 0.0001418096397599883
$ws.Range("I655").NumberFormat = "General"
$ws.Range("J655").Value = 228.5704225760452

$ws.Range("A656").Value = "'2024-09-24"
$ws.Range("C656").Value = 1904.650024414062
$ws.Range("D656").Value = 646.8499755859375
$ws.Range("E656").Value = 1051.550048828125
$ws.Range("F656").Value = 158.7400054931641
$ws.Range("G656").Value = 675.25
$ws.Range("H656").Value = 22510.13034057617
$ws.Range("I656").Value = -0.005537294613544976
$ws.Range("I656").NumberFormat = "General"
$ws.Range("J656").Value = 227.3047608062991

$ws.Range("A657").Value = "'2024-09-25"
$ws.Range("C657").Value = 1928.5
$ws.Range("D657").Value = 633.2999877929688
$ws.Range("E657").Value = 1063.449951171875
$ws.Range("F657").Value = 156.9400024414062
$ws.Range("G657").Value = 667.3499755859375
$ws.Range("H657").Value = 22551.57971191406
$ws.Range("I657").Value = 0.001841365230265907
$ws.Range("I657").NumberFormat = "General"
$ws.Range("J657").Value = 227.7233118895218
